# Add three new pin-definition rows to the "pin definition" sheet, per
# Sam's request: Temperature decode / Voltage decode / Reset Output rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pin definition")

$ws.Range("A14").Value = "Temperature decode"
$ws.Range("A15").Value = "Voltage decode"
$ws.Range("B14").Value = "RP66_RD2"
$ws.Range("B15").Value = "RP69_RD5"
$ws.Range("A16").Value = "Reset Output"
$ws.Range("B16").Value = "RP72_RD8"

# Column B carries the bold "Calibri" style used throughout the sheet.
$ws.Range("B14:B16").Font.Bold = $true

$ws.Range("B19").Select()
